$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete WELL data rows (rows 11 and 12), shifting
# everything below them up by two rows. This turns the old layout:
#   row 11: 5.8, 600, 35, 0.1, 5
#   row 12: 9.3, 2000, 40, 0.1, 5
#   row 13: RESV ...
# into:
#   row 11: RESV ...
# and shrinks the used range from A1:H19 to A1:H17.
$ws.Rows("11:12").Delete()
